$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("client_data")

# --- Insert new columns -------------------------------------------------
# Old layout: A db_id | B name_titel_vname | C name_nname | D post_str | E post_plz | F post_ort | ...
# New layout: A db_id | B titel | C vorname | D nachname | E anrede_adressfeld | F anrede | G post_str | ...
#
# Step 1: insert a single column before old B so old B/C (name_titel_vname /
#         name_nname) become new C/D (renamed to vorname / nachname), and a
#         fresh empty column B is created for "titel". All columns from D on
#         (post_str etc.) shift right along with their data/styles intact.
$ws.Range("B1").EntireColumn.Insert()

# Step 2: insert two columns before the old "post_str" column (now at E)
#         to create the new "anrede_adressfeld" (E) and "anrede" (F) columns.
#         Everything from post_str onward shifts right by two more columns.
$ws.Range("E1:F1").EntireColumn.Insert()

# --- Header row (row 1) ---------------------------------------------------
# A1, G1..Q1 already carry the correct (shifted) text from the original
# sheet; only the brand-new / renamed columns need to be written.
$ws.Range("B1").Value = "titel"
$ws.Range("C1").Value = "vorname"
$ws.Range("D1").Value = "nachname"
$ws.Range("E1").Value = "anrede_adressfeld"
$ws.Range("F1").Value = "anrede"

# --- Data rows --------------------------------------------------------
# Columns C (vorname) and D (nachname) already hold the correct values
# (shifted from the old name_titel_vname / name_nname columns), as do every
# column from G onward. Only "titel" (B) and the two new salutation columns
# (E anrede_adressfeld, F anrede) need values.

# Row 2 - client 1 (John1 Doe1, male)
$ws.Range("E2").Value = "Herrn"
$ws.Range("F2").Value = "er Herr"

# Row 3 - client 2 (Jane2 Doe2, female)
$ws.Range("E3").Value = "Frau"
$ws.Range("F3").Value = "e Frau"

# Row 4 - client 3 (John3 Doe3, male)
$ws.Range("E4").Value = "Herrn"
$ws.Range("F4").Value = "Herr"

# Row 5 - client 4 (Jane4 Doe4, female, has a title)
$ws.Range("B5").Value = "Dr."
$ws.Range("E5").Value = "Frau"
$ws.Range("F5").Value = "e Frau"

# --- window position (cosmetic, matches target OOXML) --------------------
$wb.Windows.Item(1).Left = 4820
$wb.Windows.Item(1).Top = 5020

# --- active selection moves to G10 on the client_data sheet ---------------
$ws.Range("G10").Select()
